$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1018
$ws.Range("I6").Value = 1100
$ws.Range("K6").Value = 3300
$ws.Range("M6").Value = -3188
$ws.Range("H62").Value = 6469.727
$ws.Range("I62").Value = 6777.4
$ws.Range("J62").Value = 6213.3335
$ws.Range("K62").Value = 6777.4
$ws.Range("L62").Value = 6213.3335
$ws.Range("M62").Value = -6153.4
$ws.Range("N62").Value = -7461.3335
$ws.Range("H65").Value = 6469.727
$ws.Range("I65").Value = 6777.4
$ws.Range("J65").Value = 6213.3335
$ws.Range("K65").Value = 33887
$ws.Range("L65").Value = 31066.6675
$ws.Range("M65").Value = -30767
$ws.Range("N65").Value = -37306.6675
$ws.Range("H92").Value = 99.833336
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H107").Value = 663.8946999999999
$ws.Range("I107").Value = 486
$ws.Range("J107").Value = 908.5
$ws.Range("K107").Value = 486
$ws.Range("L107").Value = 908.5
$ws.Range("M107").Value = 1434
$ws.Range("N107").Value = -4748.5
$ws.Range("H112").Value = 4172.467
$ws.Range("J112").Value = 4172.467
$ws.Range("L112").Value = 12517.401
$ws.Range("N112").Value = -14733.401
$ws.Range("H135").Value = 678.0769
$ws.Range("I135").Value = 591.6
$ws.Range("K135").Value = 5324.400000000001
$ws.Range("M135").Value = -2789.400000000001
$ws.Range("H136").Value = 49990
$ws.Range("J136").Value = 49990
$ws.Range("L136").Value = 49990
$ws.Range("N136").Value = -60190
$ws.Range("H137").Value = 2063.7942
$ws.Range("I137").Value = 1832.3914
$ws.Range("J137").Value = 2547.6365
$ws.Range("K137").Value = 5497.174199999999
$ws.Range("L137").Value = 7642.9095
$ws.Range("M137").Value = -2947.174199999999
$ws.Range("N137").Value = -12742.9095
$ws.Range("H138").Value = 3979.8635
$ws.Range("I138").Value = 10179.4
$ws.Range("J138").Value = 2156.4707
$ws.Range("K138").Value = 30538.2
$ws.Range("L138").Value = 6469.4121
$ws.Range("M138").Value = -25398.2
$ws.Range("N138").Value = -16749.4121
$ws.Range("H140").Value = 77431
$ws.Range("J140").Value = 77431
$ws.Range("L140").Value = 77431
$ws.Range("N140").Value = -87791
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3749.1875
$ws.Range("I32").Value = 3332.8774
$ws.Range("K32").Value = 3332.8774
$ws.Range("M32").Value = -3045.8774
$ws.Range("H61").Value = 2479
$ws.Range("I61").Value = 1331.1111
$ws.Range("J61").Value = 5430.7144
$ws.Range("K61").Value = 1331.1111
$ws.Range("L61").Value = 5430.7144
$ws.Range("M61").Value = -1119.1111
$ws.Range("N61").Value = -5854.7144
$ws.Range("H74").Value = 933.6923
$ws.Range("I74").Value = 595.2727
$ws.Range("J74").Value = 2795
$ws.Range("K74").Value = 595.2727
$ws.Range("L74").Value = 2795
$ws.Range("M74").Value = 278.7273
$ws.Range("N74").Value = -4543
$ws.Range("H77").Value = 933.6923
$ws.Range("I77").Value = 595.2727
$ws.Range("J77").Value = 2795
$ws.Range("K77").Value = 2976.3635
$ws.Range("L77").Value = 13975
$ws.Range("M77").Value = 1391.6365
$ws.Range("N77").Value = -22711
$ws.Range("H132").Value = 1683.7858
$ws.Range("I132").Value = 1689.3636
$ws.Range("K132").Value = 5068.0908
$ws.Range("M132").Value = -2538.0908
$ws.Range("H136").Value = 2479
$ws.Range("I136").Value = 1331.1111
$ws.Range("J136").Value = 5430.7144
$ws.Range("K136").Value = 3993.3333
$ws.Range("L136").Value = 16292.1432
$ws.Range("M136").Value = -1443.3333
$ws.Range("N136").Value = -21392.1432
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1967.8636
$ws.Range("I20").Value = 1965.4667
$ws.Range("J20").Value = 1973
$ws.Range("K20").Value = 1965.4667
$ws.Range("L20").Value = 1973
$ws.Range("M20").Value = -1718.4667
$ws.Range("N20").Value = -2467
$ws.Range("H94").Value = 2324.75
$ws.Range("I94").Value = 2266.3333
$ws.Range("K94").Value = 2266.3333
$ws.Range("M94").Value = -1815.3333
$ws.Range("H105").Value = 2164
$ws.Range("I105").Value = 2172.4783
$ws.Range("J105").Value = 2099
$ws.Range("K105").Value = 2172.4783
$ws.Range("L105").Value = 2099
$ws.Range("M105").Value = -425.4783000000002
$ws.Range("N105").Value = -5593
$ws.Range("H134").Value = 3652.6538
$ws.Range("I134").Value = 3652.6538
$ws.Range("K134").Value = 10957.9614
$ws.Range("M134").Value = -8422.9614
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1448.6786
$ws.Range("I31").Value = 901
$ws.Range("J31").Value = 1996.3572
$ws.Range("K31").Value = 901
$ws.Range("L31").Value = 1996.3572
$ws.Range("M31").Value = -606
$ws.Range("N31").Value = -2586.3572
$ws.Range("H34").Value = 1448.6786
$ws.Range("I34").Value = 901
$ws.Range("J34").Value = 1996.3572
$ws.Range("K34").Value = 901
$ws.Range("L34").Value = 1996.3572
$ws.Range("M34").Value = -699
$ws.Range("N34").Value = -2400.3572
$ws.Range("H99").Value = 668490.1
$ws.Range("I99").Value = 1430381.6
$ws.Range("J99").Value = 1835.125
$ws.Range("K99").Value = 1430381.6
$ws.Range("L99").Value = 1835.125
$ws.Range("M99").Value = -1428883.6
$ws.Range("N99").Value = -4831.125
$ws.Range("H126").Value = 668490.1
$ws.Range("I126").Value = 1430381.6
$ws.Range("J126").Value = 1835.125
$ws.Range("K126").Value = 4291144.800000001
$ws.Range("L126").Value = 5505.375
$ws.Range("M126").Value = -4288674.800000001
$ws.Range("N126").Value = -10445.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3890.75
$ws.Range("I3").Value = 3854.3333
$ws.Range("K3").Value = 11562.9999
$ws.Range("M3").Value = -11450.9999
$ws.Range("H56").Value = 68654.125
$ws.Range("I56").Value = 68654.125
$ws.Range("K56").Value = 68654.125
$ws.Range("M56").Value = -68124.125
$ws.Range("H92").Value = 879
$ws.Range("I92").Value = 200
$ws.Range("J92").Value = 954.44446
$ws.Range("K92").Value = 600
$ws.Range("L92").Value = 2863.33338
$ws.Range("M92").Value = 648
$ws.Range("N92").Value = -5359.33338
$ws.Range("H107").Value = 1468.7916
$ws.Range("J107").Value = 1768
$ws.Range("L107").Value = 5304
$ws.Range("N107").Value = -9144
$ws.Range("H137").Value = 3963.4211
$ws.Range("I137").Value = 1997
$ws.Range("K137").Value = 5991
$ws.Range("M137").Value = -891
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 101002
$ws.Range("I14").Value = 101002
$ws.Range("K14").Value = 101002
$ws.Range("M14").Value = -100834
$ws.Range("H97").Value = 1600.3043
$ws.Range("J97").Value = 2623.875
$ws.Range("L97").Value = 2623.875
$ws.Range("N97").Value = -3615.875
$ws.Range("H104").Value = 49997.5
$ws.Range("J104").Value = 49997.5
$ws.Range("L104").Value = 49997.5
$ws.Range("N104").Value = -56985.5
$ws.Range("H107").Value = 668.8570999999999
$ws.Range("I107").Value = 144.75
$ws.Range("K107").Value = 144.75
$ws.Range("M107").Value = 1775.25
$ws.Range("H122").Value = 2445.8462
$ws.Range("I122").Value = 1482
$ws.Range("J122").Value = 3988
$ws.Range("K122").Value = 4446
$ws.Range("L122").Value = 11964
$ws.Range("M122").Value = -1996
$ws.Range("N122").Value = -16864
$ws.Range("H126").Value = 3089398.5
$ws.Range("I126").Value = 3971304.2
$ws.Range("K126").Value = 11913912.6
$ws.Range("M126").Value = -11911442.6
$ws.Range("H132").Value = 12823516
$ws.Range("I132").Value = 19233274
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 57699822
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -57697292
$ws.Range("N132").Value = -17057
$ws.Range("H135").Value = 51499.15
$ws.Range("J135").Value = 51499.15
$ws.Range("L135").Value = 51499.15
$ws.Range("N135").Value = -61639.15
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3822.16
$ws.Range("I16").Value = 4931.7646
$ws.Range("K16").Value = 4931.7646
$ws.Range("M16").Value = -4761.7646
$ws.Range("H68").Value = 3635
$ws.Range("I68").Value = 3407.5
$ws.Range("K68").Value = 3407.5
$ws.Range("M68").Value = -2658.5
$ws.Range("H71").Value = 3635
$ws.Range("I71").Value = 3407.5
$ws.Range("K71").Value = 17037.5
$ws.Range("M71").Value = -13293.5
$ws.Range("H122").Value = 8899.375
$ws.Range("I122").Value = 7539.6
$ws.Range("J122").Value = 11165.667
$ws.Range("K122").Value = 22618.8
$ws.Range("L122").Value = 33497.001
$ws.Range("M122").Value = -20168.8
$ws.Range("N122").Value = -38397.001
$ws.Range("H135").Value = 32214.5
$ws.Range("J135").Value = 32214.5
$ws.Range("L135").Value = 32214.5
$ws.Range("N135").Value = -42354.5
$ws.Range("H141").Value = 31970
$ws.Range("J141").Value = 31970
$ws.Range("L141").Value = 31970
$ws.Range("N141").Value = -42330
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 44998
$ws.Range("J105").Value = 44998
$ws.Range("L105").Value = 44998
$ws.Range("N105").Value = -51986
$ws.Range("H122").Value = 269393.56
$ws.Range("I122").Value = 313962.5
$ws.Range("J122").Value = 1980
$ws.Range("K122").Value = 941887.5
$ws.Range("L122").Value = 5940
$ws.Range("M122").Value = -939437.5
$ws.Range("N122").Value = -10840
$ws.Range("H132").Value = 1091.8334
$ws.Range("I132").Value = 1091.8334
$ws.Range("K132").Value = 3275.5002
$ws.Range("M132").Value = -745.5001999999999
$ws.Range("H137").Value = 46500
$ws.Range("J137").Value = 46500
$ws.Range("L137").Value = 46500
$ws.Range("N137").Value = -56700
